# "Steel scenarios from IEA, Birat, and EUROFER"
#
# Adds C-content figures (column G) for coal, charcoal and diesel, and wires
# the CO2 (column D) figures for coal/charcoal/diesel/coke/coking coal/PCI
# coal over to a C-content-driven formula (mass CO2 per mass fuel = C
# content * 44/12). Also adds a H-content-driven "water of combustion"
# formula in column E for coke/coking coal/PCI coal, refreshes a few
# HHV/LHV (columns B/C) figures with updated source data, and adds the
# Eurofer electricity-mix-proxy LHV figure (column C). Finally leaves a
# blank spacer row (23) below the table and moves the selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - coal: new HHV/LHV, C content, and CO2 now derived from C content
$ws.Range("B3").Value = 31
$ws.Range("C3").Value = 31
$ws.Range("G3").Value = 0.81
$ws.Range("D3").Formula = "=G3*(44/12)"

# Row 4 - charcoal: add C content, CO2 now derived from C content
$ws.Range("G4").Value = 0.8
$ws.Range("D4").Formula = "=G4*(44/12)"

# Row 6 - diesel: add C content, CO2 now derived from C content
$ws.Range("G6").Value = 0.86
$ws.Range("D6").Formula = "=G6*(44/12)"

# Row 9 - Eurofer electricity mix proxy: add LHV (same as HHV)
$ws.Range("C9").Value = 1

# Row 10 - coke: align HHV with LHV, CO2 from C content, H2O from H content
$ws.Range("B10").Value = 29.01
$ws.Range("D10").Formula = "=G10*(44/12)"
$ws.Range("E10").Formula = "=18/2*H10"

# Row 11 - coking coal: new HHV/LHV, CO2 from C content, H2O from H content
$ws.Range("B11").Value = 31
$ws.Range("C11").Value = 31
$ws.Range("D11").Formula = "=G11*(44/12)"
$ws.Range("E11").Formula = "=18/2*H11"

# Row 12 - PCI coal: new HHV/LHV, CO2 from C content (was =D10), H2O from H content
$ws.Range("B12").Value = 33.37
$ws.Range("C12").Value = 33.37
$ws.Range("D12").Formula = "=G12*(44/12)"
$ws.Range("E12").Formula = "=18/2*H12"

# Blank spacer row below the table, with the row height Excel leaves behind
$ws.Rows.Item(23).RowHeight = 16.5

# Selection moves down to reflect the now-larger sheet
$ws.Range("F14").Select()
